# Modificacion: mostrar ficha del empleado y asignacion al distributivo
# The underlying data fix: NumeroPartida (column A) for every data row
# (rows 2-35 on "Hoja1") changes from 213453232 to 22323232.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:A35").Value = 22323232
